$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "513.67", "1.00", "56.557.56") are stored as text, not coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '56.557.56'
$ws.Range('E2').Value = '  -3.91%  '
$ws.Range('D3').Value = '2.360.24'
$ws.Range('E3').Value = '  -6.26%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '513.67'
$ws.Range('E5').Value = '  -3.53%  '
$ws.Range('D6').Value = '127.81'
$ws.Range('E6').Value = '  -5.29%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.49%  '
$ws.Range('D8').Value = '0.553'
$ws.Range('E8').Value = '  -2.12%  '
$ws.Range('D9').Value = '2.377.12'
$ws.Range('E9').Value = '  -5.60%  '
$ws.Range('D10').Value = '0.0959'
$ws.Range('E10').Value = '  -3.42%  '
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('D12').Value = '4.82'
$ws.Range('E12').Value = '  -8.03%  '
$ws.Range('D13').Value = '0.317'
$ws.Range('E13').Value = '  -5.05%  '
$ws.Range('D14').Value = '2.798.79'
$ws.Range('E14').Value = '  -5.38%  '
$ws.Range('D15').Value = '56.548.74'
$ws.Range('E15').Value = '  -3.98%  '
$ws.Range('D16').Value = '21.47'
$ws.Range('E16').Value = '  -4.27%  '
$ws.Range('E17').Value = '  -4.15%  '
$ws.Range('D18').Value = '2.402.82'
$ws.Range('E18').Value = '  -4.15%  '
$ws.Range('D19').Value = '10.28'
$ws.Range('E19').Value = '  -3.75%  '
$ws.Range('E20').Value = '  -4.30%  '
$ws.Range('D21').Value = '310.81'
$ws.Range('E21').Value = '  -3.67%  '
$ws.Range('D22').Value = '6.11'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = '65.46'
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.45%  '
$ws.Range('D26').Value = '0.389'
$ws.Range('E26').Value = '  -5.27%  '
$ws.Range('D27').Value = '2.464.46'
$ws.Range('E27').Value = '  -5.82%  '
$ws.Range('D28').Value = '0.154'
$ws.Range('E28').Value = '  -4.32%  '
$ws.Range('D29').Value = '7.20'
$ws.Range('E29').Value = '  -3.99%  '
$ws.Range('D30').Value = '174.82'
$ws.Range('D31').Value = '1.68'
$ws.Range('E31').Value = '  -3.88%  '
$ws.Range('D32').Value = '0.0₃0715'
$ws.Range('E32').Value = '  -6.44%  '
$ws.Range('D33').Value = '6.13'
$ws.Range('E33').Value = '  -2.73%  '
$ws.Range('E34').Value = '  -6.59%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -0.26%  '
$ws.Range('D37').Value = '17.65'
$ws.Range('E37').Value = '  -3.09%  '
$ws.Range('D38').Value = '1.19'
$ws.Range('E38').Value = '  -5.47%  '
$ws.Range('D39').Value = '3.72'
$ws.Range('E39').Value = '  -6.60%  '
$ws.Range('D40').Value = '0.810'
$ws.Range('E40').Value = '  +2.90%  '
$ws.Range('D41').Value = '35.53'
$ws.Range('E41').Value = '  -3.19%  '
$ws.Range('E42').Value = '  -6.41%  '
$ws.Range('D43').Value = '3.37'
$ws.Range('E43').Value = '  -3.43%  '
$ws.Range('D44').Value = '4.87'
$ws.Range('E44').Value = '  -4.92%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = '0.568'
$ws.Range('E45').Value = '  -4.43%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = '252.79'
$ws.Range('E46').Value = '  -9.50%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '121.03'
$ws.Range('E47').Value = '  -8.30%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '0.0906'
$ws.Range('E48').Value = '  -2.71%  '
$ws.Range('E49').Value = '  -4.23%  '
$ws.Range('E50').Value = '  -4.83%  '
$ws.Range('E51').Value = '  -6.76%  '

# Restore original (default) style on the affected columns so no residual
# number-format style change is left on the cells.
$ws.Range("D2:E51").Style = "Normal"
